# draft-gandhi-mpls-ioam-sr-06.pptx -- "Add files via upload"
#
# Renames the "Next Protocol" GACH field to "Next IP Version" throughout the
# deck:
#   - splits "...| Reserved | IOAM G-ACh..." into "...| RESVD |" + bold/blue
#     "NextIPv" + "| IOAM G-ACh..." on the diagram slides that still spelled
#     the field out as "Reserved" (slides 10, 15, 25) or already had a
#     "NexProt" run (slides 6, 30)
#   - adds "(Type TBA3)" to a trailing label on slide 25
#   - resizes/repositions the explanatory callout box on slide 6 and bumps
#     its bullet text from "Next Protocol" wording to "Next IP Version"
#     wording at a slightly larger font size

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slides 10 & 15: "   |0 0 0 1|Version| Reserved      | IOAM G-" + "ACh" +
# " (Type TBA3)        |  | "  (sz=1000)
# ---------------------------------------------------------------------
foreach ($slideIdx in 10, 15) {
    $s = $p.Slides.Item($slideIdx)
    $para = $s.Shapes.Item(4).TextFrame.TextRange.Paragraphs(6)

    $whole = $para.Characters(1, 72)
    $whole.Text = "   |0 0 0 1|Version| RESVD |NextIPv| IOAM G-ACh (Type TBA3)        |  | "

    $mark = $para.Characters(29, 7)
    $mark.Font.Bold = $true
    $mark.Font.Color.RGB = 0xC07000
}

# ---------------------------------------------------------------------
# Slide 25: "   |0 0 0 1|Version|  Reserved     | IOAM G-" + "ACh" +
# "                    |  | "  (sz=800) -- also gains "(Type TBA3)" in the
# trailing label.
# ---------------------------------------------------------------------
$s25 = $p.Slides.Item(25)
$para25 = $s25.Shapes.Item(4).TextFrame.TextRange.Paragraphs(17)

$whole25 = $para25.Characters(1, 44)
$whole25.Text = "   |0 0 0 1|Version| RESVD |NextIPv| IOAM G-"

$mark25 = $para25.Characters(29, 7)
$mark25.Font.Bold = $true
$mark25.Font.Color.RGB = 0xC07000

$tail25 = $para25.Characters(48, 25)
$tail25.Text = " (Type TBA3)        |  | "

# ---------------------------------------------------------------------
# Slide 30, "Rectangle 7": already-split "RESVD |" / "NexProt" / "IOAM G-" /
# "ACh" / "(Type TBA3)" runs get merged back down to "RESVD |NextIPv|
# IOAM G-" (plain except for the bold/blue NextIPv marker) plus "ACh" and
# "(Type TBA3)        |  | " losing their bold.
# ---------------------------------------------------------------------
$s30 = $p.Slides.Item(30)
$para30 = $s30.Shapes.Item(5).TextFrame.TextRange.Paragraphs(6)

$whole30 = $para30.Characters(1, 72)
$whole30.Text = "   |0 0 0 1|Version| RESVD |NextIPv| IOAM G-ACh (Type TBA3)        |  | "

$mark30 = $para30.Characters(29, 7)
$mark30.Font.Bold = $true
$mark30.Font.Color.RGB = 0xC07000

$ach30 = $para30.Characters(45, 3)
$ach30.Font.Bold = $false

# ---------------------------------------------------------------------
# Slide 6, "Rectangle 5": the "NexProt" run (already bold/blue/err) simply
# becomes "NextIPv".
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$para6 = $s6.Shapes.Item(4).TextFrame.TextRange.Paragraphs(4)
$mark6 = $para6.Characters(29, 7)
$mark6.Text = "NextIPv"

# ---------------------------------------------------------------------
# Slide 6, "TextBox 2": reposition/resize, bump bullets to sz=1050, and
# reword the first bullet from "Next Protocol" to "Next IP Version".
# ---------------------------------------------------------------------
$tb6 = $s6.Shapes.Item(5)
$tb6.Left = 12.0
$tb6.Top = 329.57276990551185
$tb6.Width = 282.0
$tb6.Height = 70.88551181102362

$tbtr = $tb6.TextFrame.TextRange

$b1 = $tbtr.Paragraphs(1)
$b1chars = $b1.Characters(1, 39)
$b1chars.Text = "Next IP Version is added in GACH Header "
$b1.Font.Size = 10.5

$tbtr.Paragraphs(2).Font.Size = 10.5
$tbtr.Paragraphs(3).Font.Size = 10.5
$tbtr.Paragraphs(4).Font.Size = 10.5
